$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 7 ("Experimental"): value cell B7 was blank, now holds the literal
# text "false". Plain .Value assignment would be auto-typed to a real
# Boolean by the engine (like Excel does for bare TRUE/FALSE), but the
# target workbook stores it as a shared string. Force text via the T()
# formula (text coercion) then collapse the formula to its static value
# with a values-only paste so the stored cell is a plain string, not a
# formula and not a boolean.
$cell = $ws.Range("B7")
$cell.Formula = "=T(""false"")"
$cell.Copy()
$cell.PasteSpecial(-4163)

# Row 8 ("Date"): refresh the generation timestamp.
$ws.Range("B8").Value = "2025-11-30T13:08:37+00:00"

# Row 17 ("Description"): value cell B17 was blank, now populated.
$ws.Range("B17").Value = "Directions of trends in health metrics over time"
